# Update "Number of Features" values in the Sequential Algorithm(Backward) based FS
# table on slide 9 (beta values updated for SBS).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)
$tbl = $s.Shapes.Item("Table 3").Table

# Row 2 = Decision Tree, Row 3 = Random Forest, Row 4 = SVM, Row 5 = GBM
# Column 2 = "Number of Features"
$tbl.Cell(2, 2).Shape.TextFrame.TextRange.Text = "2"
$tbl.Cell(3, 2).Shape.TextFrame.TextRange.Text = "91"
$tbl.Cell(4, 2).Shape.TextFrame.TextRange.Text = "93"
$tbl.Cell(5, 2).Shape.TextFrame.TextRange.Text = "93"
